$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Экзамен")

# Fill in new homework/practice scores for row 5 (Ивлев Андрей Сергеевич)
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 4

# Fill in new homework/practice scores for row 7 (Лемягов Александр Сергеевич)
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0

# Update selection to reflect where the user last clicked
$ws.Range("H5").Select()

$wb.Save()
